$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "痴漢|ちかん"
$ws.Range("B34").Value = "ふる|振る"
$ws.Range("B40").Value = "ためる|貯める"
$ws.Range("B42").Value = "ほめる|褒める"
$ws.Range("B45").Value = "ばかにする|馬鹿にする"
